$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D look like plain decimal numbers (e.g. "0.9961").
# Assigning such a string straight to Range.Value lets Excel auto-convert it into a
# real floating point number, which can silently drop significant trailing zeros
# (e.g. "0.2950" -> 0.295). To keep these as literal text -- exactly like the
# original inline strings -- we prefix the assignment with a quote-prefix marker
# ( ' ) which forces Excel to store the value as text, and then immediately restore
# the cell to the "Normal" style so no stray formatting/quote-prefix is left behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "30.250.93"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "1.897.05"
$ws.Range("E3").Value = "  +3.91%  "

Set-TextValue $ws.Range("D4") "0.9961"
$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "247.71"
$ws.Range("E5").Value = "  +0.43%  "

Set-TextValue $ws.Range("D6") "0.9969"
$ws.Range("E6").Value = "  -0.13%  "

Set-TextValue $ws.Range("D7") "0.4993"
$ws.Range("E7").Value = "  +1.38%  "

Set-TextValue $ws.Range("D8") "44.84"
$ws.Range("E8").Value = "  +0.58%  "

Set-TextValue $ws.Range("D9") "0.2950"
$ws.Range("E9").Value = "  +6.35%  "

Set-TextValue $ws.Range("D10") "0.06671"
$ws.Range("E10").Value = "  +4.30%  "

$ws.Range("D11").Value = "1.891.30"
$ws.Range("E11").Value = "  +3.68%  "

Set-TextValue $ws.Range("D12") "17.03"
$ws.Range("E12").Value = "  +2.01%  "

Set-TextValue $ws.Range("D13") "0.07206"
$ws.Range("E13").Value = "  +1.65%  "

Set-TextValue $ws.Range("D14") "0.6774"
$ws.Range("E14").Value = "  +4.97%  "

Set-TextValue $ws.Range("D15") "86.06"
$ws.Range("E15").Value = "  +2.19%  "

Set-TextValue $ws.Range("D16") "4.863"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").Value = "30.175.52"
$ws.Range("E17").Value = "  +3.52%  "

Set-TextValue $ws.Range("D18") "0.000008014"

Set-TextValue $ws.Range("D19") "0.9989"
$ws.Range("E19").Value = "  -0.01%  "

Set-TextValue $ws.Range("D20") "12.93"
$ws.Range("E20").Value = "  +5.83%  "

$ws.Range("D21").Value = "2.136.63"
$ws.Range("E21").Value = "  +3.88%  "

Set-TextValue $ws.Range("D22") "0.9986"
$ws.Range("E22").Value = "  +0.16%  "

Set-TextValue $ws.Range("D23") "4.781"
$ws.Range("E23").Value = "  +5.17%  "

Set-TextValue $ws.Range("D24") "5.658"
$ws.Range("E24").Value = "  +5.09%  "

Set-TextValue $ws.Range("D25") "9.173"
$ws.Range("E25").Value = "  +3.68%  "

Set-TextValue $ws.Range("D26") "147.14"
$ws.Range("E26").Value = "  +2.43%  "

Set-TextValue $ws.Range("D27") "133.81"
$ws.Range("E27").Value = "  +1.65%  "

Set-TextValue $ws.Range("D28") "16.84"
$ws.Range("E28").Value = "  +2.48%  "

Set-TextValue $ws.Range("D29") "1.952"
$ws.Range("E29").Value = "  +3.21%  "

Set-TextValue $ws.Range("D30") "1.378"
$ws.Range("E30").Value = "  -1.41%  "

Set-TextValue $ws.Range("D31") "4.232"
$ws.Range("E31").Value = "  +2.64%  "

Set-TextValue $ws.Range("D32") "0.08750"
$ws.Range("E32").Value = "  +4.66%  "

Set-TextValue $ws.Range("D33") "3.963"
$ws.Range("E33").Value = "  +4.97%  "

Set-TextValue $ws.Range("D34") "0.05161"
$ws.Range("E34").Value = "  +4.20%  "

Set-TextValue $ws.Range("D35") "1.123"
$ws.Range("E35").Value = "  +2.67%  "

Set-TextValue $ws.Range("D36") "0.7103"
$ws.Range("E36").Value = "  +5.95%  "

Set-TextValue $ws.Range("D37") "2.666"
$ws.Range("E37").Value = "  -1.15%  "

Set-TextValue $ws.Range("D38") "2.780"
$ws.Range("E38").Value = "  +3.29%  "

Set-TextValue $ws.Range("D39") "2.241"
$ws.Range("E39").Value = "  -2.56%  "

Set-TextValue $ws.Range("D40") "0.9416"
$ws.Range("E40").Value = "  -1.36%  "

$ws.Range("E41").Value = "  +4.50%  "

Set-TextValue $ws.Range("D42") "6.097"
$ws.Range("E42").Value = "  -1.21%  "

$ws.Range("E43").Value = "  -0.20%  "

Set-TextValue $ws.Range("D44") "0.4217"
$ws.Range("E44").Value = "  +3.51%  "

Set-TextValue $ws.Range("D45") "103.23"
$ws.Range("E45").Value = "  +1.42%  "

Set-TextValue $ws.Range("D46") "7.523"
$ws.Range("E46").Value = "  +4.68%  "

Set-TextValue $ws.Range("D47") "0.1270"
$ws.Range("E47").Value = "  +4.18%  "

Set-TextValue $ws.Range("D48") "0.05726"
$ws.Range("E48").Value = "  +3.25%  "

Set-TextValue $ws.Range("D49") "32.88"
$ws.Range("E49").Value = "  +4.10%  "

Set-TextValue $ws.Range("D50") "8.251"
$ws.Range("E50").Value = "  +1.85%  "

Set-TextValue $ws.Range("D51") "0.3749"
$ws.Range("E51").Value = "  +4.19%  "
